# Macroferia Regional de Talca - Berenjena
# A new weekly price-report row is inserted at row 26 (pushing every
# existing record, rows 26-76, down by one to rows 27-77), and the new
# row 26 is populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 26; Excel shifts rows
# 26..76 down to 27..77, carrying their values/formatting with them.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new record's data.
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "Macroferia Regional de Talca"
$ws.Range("C26").Value = "Maule"
$ws.Range("D26").Value = 44560
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 100112001
$ws.Range("G26").Value = "Berenjena"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 150
$ws.Range("K26").Value = 9000
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = 9000
$ws.Range("N26").Value = "`$/caja 50 unidades"
$ws.Range("O26").Value = "Región del Maule"
$ws.Range("P26").Value = 180
$ws.Range("Q26").Value = 50
$ws.Range("R26").Value = "Hortaliza"
